$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 corresponds to the "jr" instruction. Fill in the ALU/control
# truth-table values that were previously left blank (untested slt/sltu
# support in the ALU control logic).
$ws.Range("E12").Value = "0"
$ws.Range("F12").Value = "0"
$ws.Range("G12").Value = "1"
$ws.Range("H12").Value = "0"
$ws.Range("I12").Value = "0"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "X"

# Update the active selection to match where the author ended up working.
$null = $ws.Range("L12").Select()
